$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.307773574878411
$ws.Range("C2").Value = 0.3900817533207714
$ws.Range("D2").Value = 0.02428095546225251
$ws.Range("F2").Value = 0.4649055614100064
$ws.Range("G2").Value = 0.3088512092434641
$ws.Range("H2").Value = 0.4806378514355387
$ws.Range("L2").Value = 0.3004181560698385
$ws.Range("O2").Value = 1.507917572930651

$ws.Range("B3").Value = 1.163580658887838
$ws.Range("C3").Value = 0.3798827275163035
$ws.Range("D3").Value = 0.02166059791213115
$ws.Range("F3").Value = 0.4657199396035665
$ws.Range("G3").Value = 0.310876240378299
$ws.Range("H3").Value = 0.4859936546155197
$ws.Range("L3").Value = 0.2888932264966115
$ws.Range("O3").Value = 1.523138207007122

$ws.Range("B4").Value = 1.074855532874778
$ws.Range("C4").Value = 0.3737088472802839
$ws.Range("D4").Value = 0.02004353464915454
$ws.Range("F4").Value = 0.466659704421275
$ws.Range("G4").Value = 0.3124891205218034
$ws.Range("H4").Value = 0.4896011118858894
$ws.Range("L4").Value = 0.2819896235270107
$ws.Range("O4").Value = 1.533924988132881

$ws.Range("B5").Value = 1.038653866740162
$ws.Range("C5").Value = 0.3712155053138133
$ws.Range("D5").Value = 0.01938256348523737
$ws.Range("F5").Value = 0.4671530215367099
$ws.Range("G5").Value = 0.3132389749249498
$ws.Range("H5").Value = 0.4911513290938103
$ws.Range("L5").Value = 0.2792198208430534
$ws.Range("O5").Value = 1.538682290512241

$ws.Range("B6").Value = 1.032639931582196
$ws.Range("C6").Value = 0.3708028608349139
$ws.Range("D6").Value = 0.01927268998591103
$ws.Range("F6").Value = 0.4672415952077174
$ws.Range("G6").Value = 0.3133690705822971
$ws.Range("H6").Value = 0.491413580548965
$ws.Range("L6").Value = 0.2787625237693163
$ws.Range("O6").Value = 1.539494052896572

$ws.Range("B7").Value = 1.074367485475648
$ws.Range("C7").Value = 0.3736751293765792
$ws.Range("D7").Value = 0.02003462861778615
$ws.Range("F7").Value = 0.4666659109359514
$ws.Range("G7").Value = 0.3124988588446556
$ws.Range("H7").Value = 0.4896216942160478
$ws.Range("L7").Value = 0.281952092940486
$ws.Range("O7").Value = 1.533987683747881

$ws.Range("B8").Value = 1.258096666694769
$ws.Range("C8").Value = 0.3865469912847459
$ws.Range("D8").Value = 0.02337917475370688
$ws.Range("F8").Value = 0.4650949594777458
$ws.Range("G8").Value = 0.3094725559183829
$ws.Range("H8").Value = 0.4824182942353659
$ws.Range("L8").Value = 0.2964085207623839
$ws.Range("O8").Value = 1.512866027298784

$ws.Range("B9").Value = 1.616795846373179
$ws.Range("C9").Value = 0.4124762808711466
$ws.Range("D9").Value = 0.02987144223593674
$ws.Range("F9").Value = 0.4655139952445495
$ws.Range("G9").Value = 0.3064846284897342
$ws.Range("H9").Value = 0.4708260694431985
$ws.Range("L9").Value = 0.3261286962867445
$ws.Range("O9").Value = 1.482919567350265

$ws.Range("B10").Value = 1.879273396793678
$ws.Range("C10").Value = 0.4319299548987487
$ws.Range("D10").Value = 0.03459903524343133
$ws.Range("F10").Value = 0.4679709947772395
$ws.Range("G10").Value = 0.306106652337661
$ws.Range("H10").Value = 0.4638576606222102
$ws.Range("L10").Value = 0.3488034166053779
$ws.Range("O10").Value = 1.467964934652713

$ws.Range("B11").Value = 1.998434730994632
$ws.Range("C11").Value = 0.4408644126937702
$ws.Range("D11").Value = 0.03674020731091332
$ws.Range("F11").Value = 0.4695588054781297
$ws.Range("G11").Value = 0.3063337869565004
$ws.Range("H11").Value = 0.4610246035513512
$ws.Range("L11").Value = 0.3593018870141691
$ws.Range("O11").Value = 1.462703222970248

$ws.Range("B12").Value = 2.043521427352005
$ws.Range("C12").Value = 0.4442595652828345
$ws.Range("D12").Value = 0.03754961757955755
$ws.Range("F12").Value = 0.4702279196611556
$ws.Range("G12").Value = 0.306477534336409
$ws.Range("H12").Value = 0.4600003149512588
$ws.Range("L12").Value = 0.3633038037005178
$ws.Range("O12").Value = 1.460933271418156

$ws.Range("B13").Value = 2.033812890142258
$ws.Range("C13").Value = 0.4435278359674157
$ws.Range("D13").Value = 0.0373753598596096
$ws.Range("F13").Value = 0.4700807926852448
$ws.Range("G13").Value = 0.3064440023831452
$ws.Range("H13").Value = 0.4602187542070055
$ws.Range("L13").Value = 0.3624407468575583
$ws.Range("O13").Value = 1.461304549875308

$ws.Range("B14").Value = 2.002144798462439
$ws.Range("C14").Value = 0.4411434980253546
$ws.Range("D14").Value = 0.03680682641976318
$ws.Range("F14").Value = 0.4696124928089205
$ws.Range("G14").Value = 0.3063444540226214
$ws.Range("H14").Value = 0.4609393616587028
$ws.Range("L14").Value = 0.3596305987181978
$ws.Range("O14").Value = 1.462553142631663

$ws.Range("B15").Value = 1.982742271963957
$ws.Range("C15").Value = 0.4396845576891337
$ws.Range("D15").Value = 0.03645839884603674
$ws.Range("F15").Value = 0.4693344881857939
$ws.Range("G15").Value = 0.3062910070964051
$ws.Range("H15").Value = 0.4613870768093449
$ws.Range("L15").Value = 0.3579127355558001
$ws.Range("O15").Value = 1.463346950747734

$ws.Range("B16").Value = 1.871480862916485
$ws.Range("C16").Value = 0.4313477460724187
$ws.Range("D16").Value = 0.03445891062359863
$ws.Range("F16").Value = 0.4678767082945186
$ws.Range("G16").Value = 0.3060998689562524
$ws.Range("H16").Value = 0.4640495939996043
$ws.Range("L16").Value = 0.3481210092149354
$ws.Range("O16").Value = 1.468339891910119

$ws.Range("B17").Value = 1.80316210610647
$ws.Range("C17").Value = 0.4262548803253594
$ws.Range("D17").Value = 0.03322983851037975
$ws.Range("F17").Value = 0.4671029879324919
$ws.Range("G17").Value = 0.3060850747747423
$ws.Range("H17").Value = 0.4657693103968157
$ws.Range("L17").Value = 0.3421611109637439
$ws.Range("O17").Value = 1.471798299690306

$ws.Range("B18").Value = 1.76384442356516
$ws.Range("C18").Value = 0.4233336082697292
$ws.Range("D18").Value = 0.03252202367725232
$ws.Range("F18").Value = 0.4667021989525466
$ws.Range("G18").Value = 0.3061141079091954
$ws.Range("H18").Value = 0.4667901522081195
$ws.Range("L18").Value = 0.338750420788017
$ws.Range("O18").Value = 1.473932497339021

$ws.Range("B19").Value = 1.75052833568634
$ws.Range("C19").Value = 0.4223459018121503
$ws.Range("D19").Value = 0.03228221921200003
$ws.Range("F19").Value = 0.4665740885476808
$ws.Range("G19").Value = 0.3061303752769788
$ws.Range("H19").Value = 0.4671412348144557
$ws.Range("L19").Value = 0.3375985894872571
$ws.Range("O19").Value = 1.474679978037699

$ws.Range("B20").Value = 1.810437105131712
$ws.Range("C20").Value = 0.4267961980811776
$ws.Range("D20").Value = 0.03336076728876236
$ws.Range("F20").Value = 0.4671807719074366
$ws.Range("G20").Value = 0.3060827616076622
$ws.Range("H20").Value = 0.4655829615598179
$ws.Range("L20").Value = 0.3427937635351128
$ws.Range("O20").Value = 1.471415131307879

$ws.Range("B21").Value = 2.011447511380823
$ws.Range("C21").Value = 0.441843516478599
$ws.Range("D21").Value = 0.03697385699224753
$ws.Range("F21").Value = 0.4697482006190299
$ws.Range("G21").Value = 0.3063721240656321
$ws.Range("H21").Value = 0.4607263840794928
$ws.Range("L21").Value = 0.3604552918190507
$ws.Range("O21").Value = 1.462180353593851

$ws.Range("B22").Value = 2.142601974944114
$ws.Range("C22").Value = 0.4517468269014557
$ws.Range("D22").Value = 0.03932700751278162
$ws.Range("F22").Value = 0.4718217259755164
$ws.Range("G22").Value = 0.3068979317388596
$ws.Range("H22").Value = 0.4578352295499286
$ws.Range("L22").Value = 0.3721518336306957
$ws.Range("O22").Value = 1.457442473656073

$ws.Range("B23").Value = 2.072623072399097
$ws.Range("C23").Value = 0.4464550393188063
$ws.Range("D23").Value = 0.03807185420419046
$ws.Range("F23").Value = 0.4706787726612447
$ws.Range("G23").Value = 0.3065863760381404
$ws.Range("H23").Value = 0.4593523808993325
$ws.Range("L23").Value = 0.3658951148356948
$ws.Range("O23").Value = 1.459852138847737

$ws.Range("B24").Value = 1.807148203479073
$ws.Range("C24").Value = 0.4265514474761858
$ws.Range("D24").Value = 0.033301578137646
$ws.Range("F24").Value = 0.4671454686307754
$ws.Range("G24").Value = 0.3060836904872843
$ws.Range("H24").Value = 0.4656671097166409
$ws.Range("L24").Value = 0.3425076922957402
$ws.Range("O24").Value = 1.471587907292502

$ws.Range("B25").Value = 1.519937911518127
$ws.Range("C25").Value = 0.4053897558611936
$ws.Range("D25").Value = 0.0281224063817973
$ws.Range("F25").Value = 0.4650243457696064
$ws.Range("G25").Value = 0.3069753145685823
$ws.Range("H25").Value = 0.4736904338342995
$ws.Range("L25").Value = 0.3179414573047978
$ws.Range("O25").Value = 1.489787135212225
